# child_16_18.xlsx — fix the "school_Sfoour" typo (-> "school_Sfour") and
# record the analyst's last cell selection.
#
# Column C holds the age-group label for each 26-row county block
# (school_Sone, school_Stwo, school_Sthree, school_Sfoour, shool_Sfive, ...).
# The 4th label in every block was misspelled "school_Sfoour"; it lands on
# row 26 of the first block and repeats every 26 rows through row 2600.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 26; $r -le 2600; $r += 26) {
    $ws.Cells.Item($r, 3).Value = "school_Sfour"
}

# Last active selection recorded in the sheet view.
$ws.Range("H29").Select() | Out-Null
